$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.872960248610809
$ws.Range("C2").Value = 0.1963555815017912
$ws.Range("D2").Value = 0.07765566356903264
$ws.Range("E2").Value = 0.07088614255006354
$ws.Range("G2").Value = 1.786701365424705
$ws.Range("H2").Value = 1.498567188937955
$ws.Range("L2").Value = 0.1918398389711626
$ws.Range("M2").Value = 0.3582438681843385
$ws.Range("N2").Value = 1.931241967203817
$ws.Range("B3").Value = 1.762554893034746
$ws.Range("C3").Value = 0.1728416719970767
$ws.Range("D3").Value = 0.0705480675152188
$ws.Range("E3").Value = 0.07081749018032824
$ws.Range("G3").Value = 1.753861651198122
$ws.Range("H3").Value = 1.49012166938536
$ws.Range("L3").Value = 0.1894037928128114
$ws.Range("M3").Value = 0.342196095387223
$ws.Range("N3").Value = 1.950600722405447
$ws.Range("B4").Value = 1.695823267415335
$ws.Range("C4").Value = 0.1583462071039321
$ws.Range("D4").Value = 0.06622556633037391
$ws.Range("E4").Value = 0.0707766522877975
$ws.Range("G4").Value = 1.734900209384676
$ws.Range("H4").Value = 1.485763588558569
$ws.Range("L4").Value = 0.1880063448771594
$ws.Range("M4").Value = 0.3325507831124312
$ws.Range("N4").Value = 1.963164353640742
$ws.Range("B5").Value = 1.668894504944774
$ws.Range("C5").Value = 0.1524243001255741
$ws.Range("D5").Value = 0.06447441843340584
$ws.Range("E5").Value = 0.07076033887781596
$ws.Range("G5").Value = 1.727473946864194
$ws.Range("H5").Value = 1.484195114236513
$ws.Range("L5").Value = 0.1874615761136624
$ws.Range("M5").Value = 0.3286724182590746
$ws.Range("N5").Value = 1.968454065278294
$ws.Range("B6").Value = 1.664438974683151
$ws.Range("C6").Value = 0.1514400628153965
$ws.Range("D6").Value = 0.06418425972842101
$ws.Range("E6").Value = 0.07075764981543298
$ws.Range("G6").Value = 1.726258934168214
$ws.Range("H6").Value = 1.48394718402227
$ws.Range("L6").Value = 0.1873726092647416
$ws.Range("M6").Value = 0.3280315659928092
$ws.Range("N6").Value = 1.969342670265775
$ws.Range("B7").Value = 1.695459025449168
$ws.Range("C7").Value = 0.1582664029221803
$ws.Range("D7").Value = 0.06620190820466121
$ws.Range("E7").Value = 0.07077643095253183
$ws.Range("G7").Value = 1.734798840950248
$ws.Range("H7").Value = 1.485741596303086
$ws.Range("L7").Value = 0.187998897932502
$ws.Range("M7").Value = 0.3324982670516263
$ws.Range("N7").Value = 1.963235005018809
$ws.Range("B8").Value = 1.834672275008927
$ws.Range("C8").Value = 0.1882597724520281
$ws.Range("D8").Value = 0.07519623803088393
$ws.Range("E8").Value = 0.07086219687499362
$ws.Range("G8").Value = 1.775127692257485
$ws.Range("H8").Value = 1.495483053966154
$ws.Range("L8").Value = 0.1909794791072414
$ws.Range("M8").Value = 0.3526673128223052
$ws.Range("N8").Value = 1.937776071206017
$ws.Range("B9").Value = 2.116126155151392
$ws.Range("C9").Value = 0.2466361843429468
$ws.Range("D9").Value = 0.09317191867349095
$ws.Range("E9").Value = 0.07104093407083534
$ws.Range("G9").Value = 1.863832824641207
$ws.Range("H9").Value = 1.521182239692678
$ws.Range("L9").Value = 0.1976055556390293
$ws.Range("M9").Value = 0.3938790778607029
$ws.Range("N9").Value = 1.893241385440454
$ws.Range("B10").Value = 2.328180379766422
$ws.Range("C10").Value = 0.2892861966680016
$ws.Range("D10").Value = 0.1065968026941846
$ws.Range("E10").Value = 0.07117886622272618
$ws.Range("G10").Value = 1.934991037237126
$ws.Range("H10").Value = 1.544131254348315
$ws.Range("L10").Value = 0.2029525480059249
$ws.Range("M10").Value = 0.4251862119236094
$ws.Range("N10").Value = 1.863829369761532
$ws.Range("B11").Value = 2.425820301383339
$ws.Range("C11").Value = 0.3086434971294807
$ws.Range("D11").Value = 0.1127542338271894
$ws.Range("E11").Value = 0.07124309298350173
$ws.Range("G11").Value = 1.968689784792048
$ws.Range("H11").Value = 1.555465138569048
$ws.Range("L11").Value = 0.2054896377184008
$ws.Range("M11").Value = 0.4396559449965665
$ws.Range("N11").Value = 1.851172032167781
$ws.Range("B12").Value = 2.462964716851957
$ws.Range("C12").Value = 0.3159676751470784
$ws.Range("D12").Value = 0.1150933327982102
$ws.Range("E12").Value = 0.07126763008881021
$ws.Range("G12").Value = 1.981643707843745
$ws.Range("H12").Value = 1.559886369278587
$ws.Range("L12").Value = 0.2064654627941138
$ws.Range("M12").Value = 0.445168283322289
$ws.Range("N12").Value = 1.846483337284631
$ws.Range("B13").Value = 2.454957407703205
$ws.Range("C13").Value = 0.3143905437206342
$ws.Range("D13").Value = 0.1145892335671022
$ws.Range("E13").Value = 0.0712623359382647
$ws.Range("G13").Value = 1.978845237241586
$ws.Range("H13").Value = 1.558928414007227
$ws.Range("L13").Value = 0.2062546300950174
$ws.Range("M13").Value = 0.4439796335430657
$ws.Range("N13").Value = 1.847488482051247
$ws.Range("B14").Value = 2.428872774370518
$ws.Range("C14").Value = 0.3092461799830915
$ws.Range("D14").Value = 0.1129465233594686
$ws.Range("E14").Value = 0.07124510732361244
$ws.Range("G14").Value = 1.969751634096582
$ws.Range("H14").Value = 1.555826279533164
$ws.Range("L14").Value = 0.2055696169681482
$ws.Range("M14").Value = 0.4401087866577384
$ws.Range("N14").Value = 1.850784195952571
$ws.Range("B15").Value = 2.412917399362811
$ws.Range("C15").Value = 0.3060943369588358
$ws.Range("D15").Value = 0.1119412864023275
$ws.Range("E15").Value = 0.07123458249396908
$ws.Range("G15").Value = 1.964206722680188
$ws.Range("H15").Value = 1.553942997992721
$ws.Range("L15").Value = 0.2051519918673108
$ws.Range("M15").Value = 0.4377420796720557
$ws.Range("N15").Value = 1.852816521943218
$ws.Range("B16").Value = 2.321823141434436
$ws.Range("C16").Value = 0.2880202898515449
$ws.Range("D16").Value = 0.1061954282818647
$ws.Range("E16").Value = 0.07117469889450856
$ws.Range("G16").Value = 1.932815639859996
$ws.Range("H16").Value = 1.543408616256528
$ws.Range("L16").Value = 0.2027888525905723
$ws.Range("M16").Value = 0.4242451835889085
$ws.Range("N16").Value = 1.864671133756552
$ws.Range("B17").Value = 2.266241707566451
$ws.Range("C17").Value = 0.2769213385404612
$ws.Range("D17").Value = 0.1026835540840523
$ws.Range("E17").Value = 0.07113834353024107
$ws.Range("G17").Value = 1.913899836994744
$ws.Range("H17").Value = 1.537175679572755
$ws.Range("L17").Value = 0.2013659816919926
$ws.Range("M17").Value = 0.4160237903488948
$ws.Range("N17").Value = 1.872128910002701
$ws.Range("B18").Value = 2.234383291776169
$ws.Range("C18").Value = 0.2705333130896292
$ws.Range("D18").Value = 0.1006683491301317
$ws.Range("E18").Value = 0.07111757208102532
$ws.Range("G18").Value = 1.903144899113016
$ws.Range("H18").Value = 1.533674786427866
$ws.Range("L18").Value = 0.2005574387031288
$ws.Range("M18").Value = 0.4113165078095662
$ws.Range("N18").Value = 1.876486385364231
$ws.Range("B19").Value = 2.22361552918693
$ws.Range("C19").Value = 0.2683697044082578
$ws.Range("D19").Value = 0.09998684302671279
$ws.Range("E19").Value = 0.07111056304532681
$ws.Range("G19").Value = 1.899524861568437
$ws.Range("H19").Value = 1.532503870698235
$ws.Range("L19").Value = 0.2002853717211082
$ws.Range("M19").Value = 0.4097263805779008
$ws.Range("N19").Value = 1.877973408886447
$ws.Range("B20").Value = 2.272146994943625
$ws.Range("C20").Value = 0.2781032730786137
$ws.Range("D20").Value = 0.1030569081096218
$ws.Range("E20").Value = 0.07114219919391829
$ws.Range("G20").Value = 1.91590051287352
$ws.Range("H20").Value = 1.537830474154362
$ws.Range("L20").Value = 0.2015164286169266
$ws.Range("M20").Value = 0.4168967508619374
$ws.Range("N20").Value = 1.871327979873485
$ws.Range("B21").Value = 2.436529831340238
$ws.Range("C21").Value = 0.3107573644596471
$ws.Range("D21").Value = 0.1134288246174577
$ws.Range("E21").Value = 0.07125016190445699
$ws.Range("G21").Value = 1.972417394217956
$ws.Range("H21").Value = 1.556733935550568
$ws.Range("L21").Value = 0.205770412330196
$ws.Range("M21").Value = 0.4412448520789809
$ws.Range("N21").Value = 1.849813328505952
$ws.Range("B22").Value = 2.544956752231542
$ws.Range("C22").Value = 0.3320638383926848
$ws.Range("D22").Value = 0.1202507623853677
$ws.Range("E22").Value = 0.0713219811797785
$ws.Range("G22").Value = 2.010479936050842
$ws.Range("H22").Value = 1.569842654815687
$ws.Range("L22").Value = 0.2086385724181525
$ws.Range("M22").Value = 0.4573499751690093
$ws.Range("N22").Value = 1.836360808887768
$ws.Range("B23").Value = 2.486995951837116
$ws.Range("C23").Value = 0.320695237231007
$ws.Range("D23").Value = 0.1166057485665135
$ws.Range("E23").Value = 0.07128353366432716
$ws.Range("G23").Value = 1.990061624783237
$ws.Range("H23").Value = 1.562777025917683
$ws.Range("L23").Value = 0.2070997262322862
$ws.Range("M23").Value = 0.4487367164310001
$ws.Range("N23").Value = 1.843484824483163
$ws.Range("B24").Value = 2.269476915314044
$ws.Range("C24").Value = 0.2775689425875782
$ws.Range("D24").Value = 0.1028881028745587
$ws.Range("E24").Value = 0.07114045564452898
$ws.Range("G24").Value = 1.914995633650364
$ws.Range("H24").Value = 1.537534184594307
$ws.Range("L24").Value = 0.20144838202215
$ws.Range("M24").Value = 0.4165020252941005
$ws.Range("N24").Value = 1.871689862697714
$ws.Range("B25").Value = 2.03906833187159
$ws.Range("C25").Value = 0.2308879563119319
$ws.Range("D25").Value = 0.08827159958416075
$ws.Range("E25").Value = 0.07099143290356569
$ws.Range("G25").Value = 1.838793307157516
$ws.Range("H25").Value = 1.513518984871752
$ws.Range("L25").Value = 0.1957291340121898
$ws.Range("M25").Value = 0.3825507504169465
$ws.Range("N25").Value = 1.90470993485097